$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had a couple of blank/placeholder rows above and
# inside the "Vendedor/Vendas" table (row 1, row 2 and row 4 relative to
# the header that lived on row 3). Removing them shifts the header up to
# row 1 and the vendor rows directly underneath it, which also carries
# each cell's existing formatting/style up with it - matching the final
# layout where the table starts at row 1 with no gap row before "vendedor1".

# Remove the two leading blank rows (old rows 1 and 2), moving the
# "Vendedor"/"Vendas" header (old row 3) up to row 1.
$ws.Rows("1:2").Delete() | Out-Null

# Remove the blank gap row that used to sit between the header and the
# first vendor row (old row 4, now row 2 after the first delete).
$ws.Rows("2:2").Delete() | Out-Null
